$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "{'Computer Science': 92.6, 'Data Structures': 82.2, 'Algorithms': 73.5, 'Machine Learning': 61.8, 'Web Development': 86.0, 'Database Systems': 78.0}"
$ws.Range("I3").Value = "{'Computer Science': 87.8, 'Data Structures': 74.1, 'Algorithms': 92.7, 'Machine Learning': 65.2, 'Web Development': 65.4, 'Database Systems': 99.6}"
$ws.Range("I4").Value = "{'Computer Science': 60.8, 'Data Structures': 88.1, 'Algorithms': 65.1, 'Machine Learning': 83.5, 'Web Development': 70.7, 'Database Systems': 60.3}"
$ws.Range("I5").Value = "{'Computer Science': 93.7, 'Data Structures': 69.9, 'Algorithms': 71.0, 'Machine Learning': 72.0, 'Web Development': 79.6, 'Database Systems': 81.6}"
$ws.Range("I6").Value = "{'Computer Science': 72.7, 'Data Structures': 84.4, 'Algorithms': 88.9, 'Machine Learning': 75.2, 'Web Development': 92.6, 'Database Systems': 66.7}"
